$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1955719557195572
$ws.Range("C2").Value = 0.4833948339483395
$ws.Range("J2").Value = 0.01107011070110701
$ws.Range("P2").Value = 0.1845018450184502
$ws.Range("S2").Value = 0.1254612546125461
$ws.Range("B3").Value = 0.007407407407407408
$ws.Range("C3").Value = 0.02222222222222222
$ws.Range("J3").Value = 0.03703703703703703
$ws.Range("P3").Value = 0.6888888888888889
$ws.Range("S3").Value = 0.2444444444444444
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.05627705627705628
$ws.Range("D6").Value = 0.008658008658008658
$ws.Range("F6").Value = 0.05194805194805195
$ws.Range("J6").Value = 0.2337662337662338
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.1341991341991342
$ws.Range("R6").Value = 0.07792207792207792
$ws.Range("S6").Value = 0.4285714285714285
$ws.Range("B7").Value = 0.08947368421052632
$ws.Range("D7").Value = 0.01578947368421053
$ws.Range("F7").Value = 0.06842105263157895
$ws.Range("J7").Value = 0.1421052631578947
$ws.Range("O7").Value = 0.03157894736842105
$ws.Range("Q7").Value = 0.1842105263157895
$ws.Range("R7").Value = 0.05789473684210526
$ws.Range("S7").Value = 0.4105263157894737
$ws.Range("B8").Value = 0.06167400881057269
$ws.Range("D8").Value = 0.01101321585903084
$ws.Range("E8").Value = 0.002202643171806168
$ws.Range("F8").Value = 0.06387665198237885
$ws.Range("J8").Value = 0.1101321585903084
$ws.Range("O8").Value = 0.02202643171806168
$ws.Range("Q8").Value = 0.2092511013215859
$ws.Range("R8").Value = 0.09030837004405286
$ws.Range("S8").Value = 0.4295154185022027
$ws.Range("B9").Value = 0.09210526315789473
$ws.Range("D9").Value = 0.0131578947368421
$ws.Range("F9").Value = 0.06578947368421052
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("O9").Value = 0.03070175438596491
$ws.Range("Q9").Value = 0.1842105263157895
$ws.Range("R9").Value = 0.06140350877192982
$ws.Range("S9").Value = 0.4473684210526316
$ws.Range("B10").Value = 0.1039729501267963
$ws.Range("D10").Value = 0.02282333051563821
$ws.Range("F10").Value = 0.08368554522400676
$ws.Range("J10").Value = 0.1191885038038884
$ws.Range("O10").Value = 0.01183431952662722
$ws.Range("Q10").Value = 0.2071005917159763
$ws.Range("R10").Value = 0.07523245984784446
$ws.Range("S10").Value = 0.3761622992392223
$ws.Range("G11").Value = 0.1338028169014084
$ws.Range("J11").Value = 0.07746478873239436
$ws.Range("K11").Value = 0.1584507042253521
$ws.Range("L11").Value = 0.6056338028169014
$ws.Range("S11").Value = 0.02464788732394366
$ws.Range("G12").Value = 0.7444444444444445
$ws.Range("J12").Value = 0.1555555555555556
$ws.Range("K12").Value = 0.01111111111111111
$ws.Range("L12").Value = 0.03333333333333333
$ws.Range("S12").Value = 0.05555555555555555
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6046511627906976
$ws.Range("J13").Value = 0.3255813953488372
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.0128755364806867
$ws.Range("H15").Value = 0.1459227467811159
$ws.Range("I15").Value = 0.07296137339055794
$ws.Range("J15").Value = 0.369098712446352
$ws.Range("K15").Value = 0.05579399141630902
$ws.Range("M15").Value = 0.01716738197424893
$ws.Range("O15").Value = 0.06008583690987124
$ws.Range("S15").Value = 0.2660944206008584
$ws.Range("F16").Value = 0.05095541401273886
$ws.Range("H16").Value = 0.1847133757961783
$ws.Range("I16").Value = 0.1019108280254777
$ws.Range("J16").Value = 0.3248407643312102
$ws.Range("K16").Value = 0.1210191082802548
$ws.Range("M16").Value = 0.01910828025477707
$ws.Range("O16").Value = 0.03821656050955414
$ws.Range("S16").Value = 0.1592356687898089
$ws.Range("F17").Value = 0.01128668171557562
$ws.Range("H17").Value = 0.2054176072234763
$ws.Range("I17").Value = 0.09706546275395034
$ws.Range("J17").Value = 0.3724604966139955
$ws.Range("K17").Value = 0.108352144469526
$ws.Range("M17").Value = 0.01354401805869074
$ws.Range("N17").Value = 0.002257336343115124
$ws.Range("O17").Value = 0.0564334085778781
$ws.Range("S17").Value = 0.1331828442437923
$ws.Range("F18").Value = 0.005882352941176471
$ws.Range("H18").Value = 0.1882352941176471
$ws.Range("I18").Value = 0.09411764705882353
$ws.Range("J18").Value = 0.4411764705882353
$ws.Range("K18").Value = 0.06470588235294118
$ws.Range("M18").Value = 0.02352941176470588
$ws.Range("O18").Value = 0.08823529411764706
$ws.Range("S18").Value = 0.09411764705882353
$ws.Range("F19").Value = 0.01518602885345482
$ws.Range("H19").Value = 0.2088078967350038
$ws.Range("I19").Value = 0.1032649962034928
$ws.Range("J19").Value = 0.3439635535307517
$ws.Range("K19").Value = 0.105542900531511
$ws.Range("M19").Value = 0.01974183750949127
$ws.Range("O19").Value = 0.08428246013667426
$ws.Range("S19").Value = 0.1192103264996204
